# Updated cryptos list on Thu Nov 21 19:32:21 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rng, $value) {
    # Force the cell to stay a text value even when it looks like a plain
    # number (e.g. "255.79"), without permanently altering the cell style.
    $ws.Range($rng).Value = "'" + $value
    $ws.Range($rng).Style = "Normal"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "98.701.88"
$ws.Range("E2").Value = "  +5.46%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.349.36"
$ws.Range("E3").Value = "  +10.01%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.09%  "

# Row 5 - Solana
Set-TextValue "D5" "255.79"
$ws.Range("E5").Value = "  +10.05%  "

# Row 6 - BNB
Set-TextValue "D6" "622.27"
$ws.Range("E6").Value = "  +3.25%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  +11.58%  "

# Row 8 - Dogecoin
Set-TextValue "D8" "0.385"
$ws.Range("E8").Value = "  +4.09%  "

# Row 9 - USDC
Set-TextValue "D9" "0.999"

# Row 10 - LidoStakedEther
$ws.Range("D10").Value = "3.345.16"
$ws.Range("E10").Value = "  +10.06%  "

# Row 11 - Cardano
Set-TextValue "D11" "0.802"
$ws.Range("E11").Value = "  +1.03%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  +2.81%  "

# Row 13 - WrappedBTC
$ws.Range("D13").Value = "98.321.73"
$ws.Range("E13").Value = "  +5.51%  "

# Row 14 - Avalanche
Set-TextValue "D14" "35.78"
$ws.Range("E14").Value = "  +8.34%  "

# Row 15 & 16 - swap WrappedliquidstakedEther2.0 and ShibaInu
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D15" "0.0000245"
$ws.Range("E15").Value = "  +3.90%  "

$ws.Range("B16").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C16").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D16").Value = "3.970.67"
$ws.Range("E16").Value = "  +10.06%  "

# Row 17 - Toncoin
$ws.Range("E17").Value = "  +4.88%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "3.354.11"
$ws.Range("E18").Value = "  +9.99%  "

# Row 19 - SuiNetwork
Set-TextValue "D19" "3.56"
$ws.Range("E19").Value = "  +3.16%  "

# Row 20 - Chainlink
Set-TextValue "D20" "14.94"
$ws.Range("E20").Value = "  +5.55%  "

# Row 21 - BitcoinCash
Set-TextValue "D21" "485.34"
$ws.Range("E21").Value = "  +12.28%  "

# Row 22 - Polkadot
Set-TextValue "D22" "5.83"
$ws.Range("E22").Value = "  +3.92%  "

# Row 23 - PEPE
$ws.Range("E23").Value = "  +10.85%  "

# Row 24 - Uniswap
Set-TextValue "D24" "9.14"
$ws.Range("E24").Value = "  +5.52%  "

# Row 25 - NEARProtocol
Set-TextValue "D25" "5.64"
$ws.Range("E25").Value = "  +3.93%  "

# Row 26 - Litecoin
Set-TextValue "D26" "88.37"
$ws.Range("E26").Value = "  +5.48%  "

# Row 27 - Aptos
Set-TextValue "D27" "11.93"
$ws.Range("E27").Value = "  +3.32%  "

# Row 28 - WrappedeETH
$ws.Range("D28").Value = "3.529.65"
$ws.Range("E28").Value = "  +10.42%  "

# Row 30 - Cronos
$ws.Range("E30").Value = "  +7.61%  "

# Row 31 - Stellar
Set-TextValue "D31" "0.250"
$ws.Range("E31").Value = "  +2.73%  "

# Row 32 - Hedera
Set-TextValue "D32" "0.125"
$ws.Range("E32").Value = "  +4.16%  "

# Row 33 - Binance-PegBSC-USD
Set-TextValue "D33" "0.999"
$ws.Range("E33").Value = "  -10.82%  "

# Row 34 - InternetComputer(DFINITY)
Set-TextValue "D34" "9.26"
$ws.Range("E34").Value = "  +4.24%  "

# Row 35 - EthereumClassic
Set-TextValue "D35" "27.21"
$ws.Range("E35").Value = "  +8.99%  "

# Row 36 - RenderToken
Set-TextValue "D36" "7.38"
$ws.Range("E36").Value = "  -2.06%  "

# Row 37 - Bittensor
Set-TextValue "D37" "516.46"
$ws.Range("E37").Value = "  +13.41%  "

# Row 38 - Kaspa
$ws.Range("E38").Value = "  +0.05%  "

# Row 39 - PancakeSwap
$ws.Range("E39").Value = "  +3.64%  "

# Row 40 - WhiteBITCoin
Set-TextValue "D40" "24.89"
$ws.Range("E40").Value = "  +3.88%  "

# Row 41 - PolygonEcosystemToken
Set-TextValue "D41" "0.445"
$ws.Range("E41").Value = "  +4.47%  "

# Row 42 - Fetch.AI
Set-TextValue "D42" "1.25"
$ws.Range("E42").Value = "  +2.44%  "

# Row 43 - MantraDAO
Set-TextValue "D43" "3.60"
$ws.Range("E43").Value = "  -4.45%  "

# Row 44 - dogwifhat
Set-TextValue "D44" "3.23"
$ws.Range("E44").Value = "  +5.22%  "

# Row 45 - USDe
$ws.Range("E45").Value = "  +0.03%  "

# Row 46 - ARBITRUM
Set-TextValue "D46" "0.771"
$ws.Range("E46").Value = "  +18.15%  "

# Row 47 - Monero
Set-TextValue "D47" "161.27"
$ws.Range("E47").Value = "  +0.29%  "

# Row 48 - Stacks
Set-TextValue "D48" "1.90"
$ws.Range("E48").Value = "  +6.32%  "

# Row 49 - OKB
Set-TextValue "D49" "45.48"
$ws.Range("E49").Value = "  +4.37%  "

# Row 50 - ImmutableX
$ws.Range("E50").Value = "  +8.49%  "

# Row 51 - Filecoin
$ws.Range("E51").Value = "  +8.10%  "
